$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$LF = [char]10

# ---------------------------------------------------------------------------
# 1. Update row 17 ([ (left square bracket)) B/C contents and row height.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "[NOTHING WORKED]"
$ws.Range("C17").Value = "2 (escape with grave accent)"
$ws.Rows.Item(17).RowHeight = 66

# ---------------------------------------------------------------------------
# 2. Add a new "Note" column (E) with data for every existing row.
# ---------------------------------------------------------------------------

# Header E1 - copy formatting from the bold header style used by B1..D1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = "Note"

# Plain bordered cells E2:E16 and E18:E26 - copy formatting from A2 (plain bordered style)
$ws.Range("A2").Copy()
$ws.Range("E2:E16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E18:E26").PasteSpecial(-4122)  # xlPasteFormats

# E17 gets the note text with a wrap-text bordered style (same style family as column D)
$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E17").Value = "To get ACL, escape with grave accent, then do Get-Item. Then, when the item is returned, access its method GetAccessControl()"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Remove the bold-ish duplicate style from A25/A26 ("curly quotes" rows) by
#    re-applying the same wrap-text bordered format already used elsewhere
#    (column D's plain wrap-text style), which drops the redundant style and
#    keeps the content intact.
# ---------------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Column E width + freeze panes at B2 (freeze header row & first column).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 32.57

$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
